# Apply fixes to trading_book workbook and add a new trade row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo / data errors in existing rows ---

# F3: "SAR EMA MACD" -> "SAR EMA MAC"
$ws.Range("F3").Value = "SAR EMA MAC"

# B8: "USDTIRT" -> "USDT/IRT"
$ws.Range("B8").Value = "USDT/IRT"

# B9: "USDTIRT" -> "USDT/IRT"
$ws.Range("B9").Value = "USDT/IRT"

# F9: was blank -> "nan"
$ws.Range("F9").Value = "nan"

# --- Add new trade row 10 ---
$ws.Range("A10").Value = 45832
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat

$ws.Range("B10").Value = "AAVE/USDT"
$ws.Range("C10").Value = "Buy"
$ws.Range("D10").Value = 0.02
$ws.Range("E10").Value = 261.18
$ws.Range("F10").Value = "T"
$ws.Range("G10").Value = 5.2236
